$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "34.539.10"
$ws.Range("E2").Value = "  +1.69%  "
$ws.Range("D3").Value = "1.841.73"
$ws.Range("E3").Value = "  +3.80%  "
$ws.Range("E4").Value = "  -0.20%  "
$ws.Range("D5").Value = "227.01"
$ws.Range("E5").Value = "  +0.81%  "
$ws.Range("D6").Value = "0.555"
$ws.Range("E6").Value = "  +1.43%  "
$ws.Range("E7").Value = "  -0.18%  "
$ws.Range("D8").Value = "32.48"
$ws.Range("E8").Value = "  +3.70%  "
$ws.Range("D9").Value = "0.296"
$ws.Range("E9").Value = "  +5.88%  "
$ws.Range("D10").Value = "0.0722"
$ws.Range("E10").Value = "  +10.29%  "
$ws.Range("D11").Value = "0.0933"
$ws.Range("E11").Value = "  +0.58%  "
$ws.Range("D12").Value = "2.105.32"
$ws.Range("E12").Value = "  +3.78%  "
$ws.Range("D13").Value = "1.848.03"
$ws.Range("E13").Value = "  +4.18%  "
$ws.Range("D14").Value = "10.99"
$ws.Range("E14").Value = "  +1.39%  "
$ws.Range("E15").Value = "  +4.87%  "
$ws.Range("D16").Value = "34.553.55"
$ws.Range("E16").Value = "  +1.66%  "
$ws.Range("D17").Value = "4.38"
$ws.Range("E17").Value = "  +4.55%  "
$ws.Range("D18").Value = "69.98"
$ws.Range("D19").Value = "253.50"
$ws.Range("E19").Value = "  +0.69%  "
$ws.Range("D20").Value = "0.0₃0809"
$ws.Range("E20").Value = "  +10.05%  "
$ws.Range("D21").Value = "11.29"
$ws.Range("E21").Value = "  +9.56%  "
$ws.Range("D22").Value = "0.998"
$ws.Range("E22").Value = "  -0.24%  "
$ws.Range("D23").Value = "4.32"
$ws.Range("E23").Value = "  +3.19%  "
$ws.Range("D24").Value = "2.16"
$ws.Range("E24").Value = "  +1.63%  "
$ws.Range("D25").Value = "161.97"
$ws.Range("D26").Value = "16.88"
$ws.Range("E26").Value = "  +3.31%  "
$ws.Range("D27").Value = "7.28"
$ws.Range("E27").Value = "  +4.27%  "
$ws.Range("E28").Value = "  +1.83%  "
$ws.Range("E29").Value = "  -0.20%  "
$ws.Range("D30").Value = "0.0539"
$ws.Range("E30").Value = "  +5.91%  "
$ws.Range("D31").Value = "3.83"
$ws.Range("E31").Value = "  +2.18%  "
$ws.Range("B32").Value = "PancakeSwap"
$ws.Range("C32").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D32").Value = "1.22"
$ws.Range("E32").Value = "  +2.36%  "
$ws.Range("B33").Value = "Swop.fi"
$ws.Range("C33").Value = "https://coinranking.com/coin/yrCr2HW2c+swopfi-swop"
$ws.Range("D33").Value = "522.48"
$ws.Range("E33").Value = "  +901.70%  "
$ws.Range("D34").Value = "3.64"
$ws.Range("E35").Value = "  +7.05%  "
$ws.Range("D36").Value = "1.466.19"
$ws.Range("E36").Value = "  +1.46%  "
$ws.Range("D37").Value = "0.657"
$ws.Range("E37").Value = "  +5.24%  "
$ws.Range("D38").Value = "1.09"
$ws.Range("E38").Value = "  +2.99%  "
$ws.Range("D39").Value = "0.0195"
$ws.Range("E39").Value = "  +4.61%  "
$ws.Range("D40").Value = "0.985"
$ws.Range("E40").Value = "  +11.33%  "
$ws.Range("D41").Value = "83.08"
$ws.Range("E41").Value = "  +0.71%  "
$ws.Range("D42").Value = "2.79"
$ws.Range("E42").Value = "  -1.45%  "
$ws.Range("E43").Value = "  +1.01%  "
$ws.Range("D44").Value = "2.16"
$ws.Range("E44").Value = "  +6.22%  "
$ws.Range("D45").Value = "6.16"
$ws.Range("E45").Value = "  +7.26%  "
$ws.Range("D46").Value = "2.002.45"
$ws.Range("E47").Value = "  +1.49%  "
$ws.Range("E48").Value = "  -1.88%  "
$ws.Range("D49").Value = "12.42"
$ws.Range("E49").Value = "  +5.13%  "
$ws.Range("D50").Value = "106.62"
$ws.Range("E50").Value = "  +9.56%  "
$ws.Range("E51").Value = "  +0.23%  "
